$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "valid_login"

# Set up header row and data row
$ws.Range("A1").Value = "username"
$ws.Range("B1").Value = "password"
$ws.Range("A2").Value = "admin"
$ws.Range("B2").Value = "manager"

# Update selection to match target (A1:B1)
$ws.Range("A1:B1").Select()

# Maximize the window (matches the large xWindow/yWindow/width/height saved in the target file)
$excel.ActiveWindow.WindowState = -4137
